$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Puan" column values (D2:D5)
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 4
$ws.Range("D4").Value = 5
$ws.Range("D5").Value = 6

# Update the active selection/cell to E3
$ws.Range("E3").Select()
